$d = $word.ActiveDocument

function Replace-InRange($range, $findText, $replaceText) {
    $r = $range.Duplicate
    $r.Find.ClearFormatting()
    [void]$r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

function Split-Sentence($range, $findFullText, $firstPart, $secondPart) {
    # Replace the whole original sentence text with the first (translated) sentence,
    # then insert a separate run containing just a space, then a separate run with
    # the second (translated) sentence - producing 3 distinct runs.
    $r = $range.Duplicate
    $r.Find.ClearFormatting()
    [void]$r.Find.Execute($findFullText, $true, $false, $false, $false, $false, $true, 1, $false, $firstPart, 2)
    $r.Collapse(0)
    $r.InsertAfter(" ")
    $r.Collapse(0)
    $r.InsertAfter($secondPart)
}

# --- Paragraph 1: title ---
Replace-InRange $d.Paragraphs.Item(1).Range "ContosoLearn Competitor SWOT" "ContosoLearn Mitbewerber SWOT"

# --- Paragraph 2: "Fabrikam Learning:" (unchanged) ---

# --- Paragraph 3: Strengths (Fabrikam) ---
$p = $d.Paragraphs.Item(3).Range
Replace-InRange $p "Strengths:" "Stärken:"
$p = $d.Paragraphs.Item(3).Range
Split-Sentence $p " Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed." " Fabrikam Learning bietet ein umfassendes Set an Analyse- und Berichterstattungs-Tools." "Es gewährleistet die kontinuierliche Überwachung von Lehr- und Lernaktivitäten sowie die Anheftung problematischer Bereiche, die angegangen werden müssen."

# --- Paragraph 4: Weaknesses (Fabrikam) ---
$p = $d.Paragraphs.Item(4).Range
Replace-InRange $p "Weaknesses:" "Schwachstellen:"
$p = $d.Paragraphs.Item(4).Range
Replace-InRange $p " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." " Während Fabrikam Learning robuste Berichterstellungsfunktionen bietet, kann es aufgrund seiner umfassenden Natur für einige Benutzende überwältigend sein."

# --- Paragraph 5: Opportunities (Fabrikam) ---
$p = $d.Paragraphs.Item(5).Range
Replace-InRange $p "Opportunities:" "Verkaufschancen:"
$p = $d.Paragraphs.Item(5).Range
Split-Sentence $p " There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand." " Es besteht eine wachsende Nachfrage nach personalisierten Lernerfahrungen und datengesteuerten Empfehlungen." "Fabrikam Learning kann seine robusten Analyse- und Berichterstellungstools nutzen, um diese Nachfrage zu erfüllen."

# --- Paragraph 6: Threats (Fabrikam) ---
$p = $d.Paragraphs.Item(6).Range
Replace-InRange $p "Threats:" "Bedrohungen:"
$p = $d.Paragraphs.Item(6).Range
Split-Sentence $p " The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead." " Der eLearning-Markt ist mit vielen Anbietern, die ähnliche Funktionen anbieten, sehr wettbewerbsfähig." "Fabrikam Learning muss kontinuierlich innovativ sein, um an der Spitze zu bleiben."

# --- Paragraph 7: "AdatumLearn:" (unchanged) ---

# --- Paragraph 8: Strengths (Adatum) ---
$p = $d.Paragraphs.Item(8).Range
Replace-InRange $p "Strengths:" "Stärken:"
$p = $d.Paragraphs.Item(8).Range
Split-Sentence $p " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users." " AdatumLearn bietet Kurse zu Geschäftsanalysetechniken wie MOST und SWOT an." "Dies zeigt ihr Engagement für die Bereitstellung wertvoller Inhalte für ihre Benutzenden."

# --- Paragraph 9: Weaknesses (Adatum) ---
$p = $d.Paragraphs.Item(9).Range
Replace-InRange $p "Weaknesses:" "Schwachstellen:"
$p = $d.Paragraphs.Item(9).Range
Split-Sentence $p " The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content." " Die in ihren Kursen bereitgestellten Informationen sind eine Zusammenstellung von generierten Informationen von Drittanbietern." "Dies ist möglicherweise nicht so wertvoll wie der ursprüngliche Inhalt."

# --- Paragraph 10: Opportunities (Adatum) ---
$p = $d.Paragraphs.Item(10).Range
Replace-InRange $p "Opportunities:" "Verkaufschancen:"
$p = $d.Paragraphs.Item(10).Range
Split-Sentence $p " AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics." " AdatumLearn kann originellere Inhalte erstellen, um ihren Benutzenden einen einzigartigen Wert zu bieten." "Sie können auch ihre Kursangebote erweitern, um weitere Themen zu behandeln."

# --- Paragraph 11: Threats (Adatum) ---
$p = $d.Paragraphs.Item(11).Range
Replace-InRange $p "Threats:" "Bedrohungen:"
$p = $d.Paragraphs.Item(11).Range
Split-Sentence $p " Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive.""" " Wie Fabrikam Learning steht AdatumLearn auch im eLearning-Markt vor einem harten Wettbewerb." "Sie müssen ihr Angebot kontinuierlich verbessern, um wettbewerbsfähig zu bleiben."

Write-Host "Done"
